$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.210.30"

$ws.Range("D3").Value = "3.774.68"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").Value = "3.772.20"
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("E13").Value = "  -1.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").Value = "4.404.14"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "3.733.02"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "68.157.33"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("E24").Value = "  +9.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").Value = "3.727.55"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.35%  "

$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("E48").Value = "  -2.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "389.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").Value = "2.785.50"
$ws.Range("E51").Value = "  +4.18%  "
